# The document contains several text boxes whose visible (DrawingML / wps)
# content lives in an <mc:Choice> branch and is mirrored by a legacy VML
# <mc:Fallback> branch. Word keeps both branches in sync; to reproduce that
# here we round-trip the whole package OOXML (Range.WordOpenXML ->
# Range.InsertXML) so every occurrence of the old text gets updated.

$d = $word.ActiveDocument

$xml = $d.Content.WordOpenXML

$xml = $xml.Replace("<w:t>scasc</w:t>", "<w:t>scvss</w:t>")
$xml = $xml.Replace("<w:t>h ji</w:t>", "<w:t>h</w:t>")
$xml = $xml.Replace("<w:t>j</w:t>", "<w:t>kjh</w:t>")
$xml = $xml.Replace("<w:t>hi ono</w:t>", "<w:t>ho</w:t>")
$xml = $xml.Replace("<w:t>iunoio</w:t>", "<w:t>uob</w:t>")
$xml = $xml.Replace("<w:t>male</w:t>", "<w:t>female</w:t>")
$xml = $xml.Replace("<w:t>2024-03-29</w:t>", "<w:t>2024-04-01</w:t>")

$d.Content.InsertXML($xml)
